$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Extend the "inventory_test" table with a new "Inventory" column
$lo = $ws1.ListObjects.Item(1)
$newCol = $lo.ListColumns.Add()
$ws1.Range("E1").Value = "Inventory"

# Populate the new Inventory column with stock values (rows 2-44)
$inventory = @(433,254,294,350,151,250,254,215,151,292,140,111,216,420,118,222,61,319,390,225,17,403,323,402,112,182,93,93,414,425,132,338,360,403,383,26,194,178,155,148,0,209,265)
for ($i = 0; $i -lt $inventory.Length; $i++) {
    $row = 2 + $i
    $ws1.Cells.Item($row, 5).Value = $inventory[$i]
}

# Give the new column a sensible width (not best-fit, matches a manually sized column)
$ws1.Columns.Item(5).ColumnWidth = 10.7109375

# Make inventory_test the active sheet/selection, matching the edited workbook
$ws1.Activate() | Out-Null
$ws1.Range("G7").Select() | Out-Null
